# Import / Export XL #25 - progress
# Rebuild the "Aclass" sheet (sheet1) with the new header row, an inserted
# CreatedAt date column, an inserted Status/enum column and a trailing
# Duration column, plus a third data row.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Aclass")

# Start from a clean sheet so stale cells from the old layout don't linger.
$ws1.Cells.Clear()

# -- Row 1: new header / country row -------------------------------------
$ws1.Range("A1").Value = "India"
$ws1.Range("B1").Value = "Canada"
$ws1.Range("C1").Value = "Japan"

# -- Row 2: first data row --------------------------------------------------
$ws1.Range("A2").Value = "A1"

$ws1.Range("B2").Value2 = 44247.84171296297
$ws1.Range("B2").NumberFormat = "m/d/yy h:mm"

$ws1.Range("C2").Value = $true
$ws1.Range("D2").Value = "ENUM_VAL1_NOT_THE_SAME"
$ws1.Range("E2").Value = ""
$ws1.Range("F2").Value = ""
$ws1.Range("G2").Value = ""
$ws1.Range("H2").Value = 0
$ws1.Range("I2").Value = 10.2
$ws1.Range("J2").Value = 4
$ws1.Range("K2").Value = $true
$ws1.Range("L2").Value = "1h3m0.001s"

# -- Row 3: second data row --------------------------------------------------
$ws1.Range("A3").Value = "A2"

$ws1.Range("B3").Value2 = -693593
$ws1.Range("B3").NumberFormat = "m/d/yy h:mm"

$ws1.Range("C3").Value = $true
$ws1.Range("D3").Value = ""
$ws1.Range("E3").Value = ""
$ws1.Range("F3").Value = ""
$ws1.Range("G3").Value = ""
$ws1.Range("H3").Value = 0
$ws1.Range("I3").Value = 10.77
$ws1.Range("J3").Value = 0
$ws1.Range("K3").Value = $true
$ws1.Range("L3").Value = "0s"

# Bclass / Dclass sheets keep the same visible values; only the shared
# string table indices move around, which the writer manages on its own.
